$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite rows 2-26 with the refreshed hotel/listing data.
# (Columns C/D/E are left blank for listings that have no review-score data yet,
#  matching the source extract.)

$ws.Cells.Item(2, 1).Value = 'Rent a Room - Charon 3'
$ws.Cells.Item(2, 2).Value = 'US$12,693'
$ws.Cells.Item(2, 3).Value = '''6.7'
$ws.Cells.Item(2, 4).Value = 'Review score'
$ws.Cells.Item(2, 5).Value = '''6'

$ws.Cells.Item(3, 1).Value = 'GuestReady - Charm and Confort in the 18th'
$ws.Cells.Item(3, 2).Value = 'US$2,443'
$ws.Cells.Item(3, 3).ClearContents()
$ws.Cells.Item(3, 4).ClearContents()
$ws.Cells.Item(3, 5).ClearContents()

$ws.Cells.Item(4, 1).Value = 'Appartment next to l''arc'
$ws.Cells.Item(4, 2).Value = 'US$2,942'
$ws.Cells.Item(4, 3).ClearContents()
$ws.Cells.Item(4, 4).ClearContents()
$ws.Cells.Item(4, 5).ClearContents()

$ws.Cells.Item(5, 1).Value = 'GuestReady - Bright Parisian Apt near Notre Dame'
$ws.Cells.Item(5, 2).Value = 'US$2,507'
$ws.Cells.Item(5, 3).ClearContents()
$ws.Cells.Item(5, 4).ClearContents()
$ws.Cells.Item(5, 5).ClearContents()

$ws.Cells.Item(6, 1).Value = 'ZZZ poteau de Paris 4P'
$ws.Cells.Item(6, 2).Value = 'US$5,678'
$ws.Cells.Item(6, 3).Value = '''7.3'
$ws.Cells.Item(6, 4).Value = 'Very Good'
$ws.Cells.Item(6, 5).Value = '''27'

$ws.Cells.Item(7, 1).Value = 'Maison le Bac Paris Aparthotel'
$ws.Cells.Item(7, 2).Value = 'US$3,735'
$ws.Cells.Item(7, 3).Value = '''8.2'
$ws.Cells.Item(7, 4).Value = 'Very Good'
$ws.Cells.Item(7, 5).Value = '''612'

$ws.Cells.Item(8, 1).Value = 'Magnifique appartement moulin rouge'
$ws.Cells.Item(8, 2).Value = 'US$2,452'
$ws.Cells.Item(8, 3).ClearContents()
$ws.Cells.Item(8, 4).ClearContents()
$ws.Cells.Item(8, 5).ClearContents()

$ws.Cells.Item(9, 1).Value = 'The Bastille Day apartment'
$ws.Cells.Item(9, 2).Value = 'US$2,564'
$ws.Cells.Item(9, 3).Value = '''7.7'
$ws.Cells.Item(9, 4).Value = 'Good'
$ws.Cells.Item(9, 5).Value = '''3'

$ws.Cells.Item(10, 1).Value = 'Appartement République Paris 11e'
$ws.Cells.Item(10, 2).Value = 'US$2,607'
$ws.Cells.Item(10, 3).Value = '''8.4'
$ws.Cells.Item(10, 4).Value = 'Very Good'
$ws.Cells.Item(10, 5).Value = '''7'

$ws.Cells.Item(11, 1).Value = 'Cozy studio in the heart of the 9th arrondissement'
$ws.Cells.Item(11, 2).Value = 'US$3,976'
$ws.Cells.Item(11, 3).ClearContents()
$ws.Cells.Item(11, 4).ClearContents()
$ws.Cells.Item(11, 5).ClearContents()

$ws.Cells.Item(12, 1).Value = 'amazing place 4 person Paris 8'
$ws.Cells.Item(12, 2).Value = 'US$15,892'
$ws.Cells.Item(12, 3).ClearContents()
$ws.Cells.Item(12, 4).ClearContents()
$ws.Cells.Item(12, 5).ClearContents()

$ws.Cells.Item(13, 1).Value = 'Cute studio in Bastille- Bail mobilité'
$ws.Cells.Item(13, 2).Value = 'US$2,446'
$ws.Cells.Item(13, 3).Value = '''8.0'
$ws.Cells.Item(13, 4).Value = 'Very Good'
$ws.Cells.Item(13, 5).Value = '''1'

$ws.Cells.Item(14, 1).Value = 'GuestReady - Stylish Hideaway near Bercy Parc'
$ws.Cells.Item(14, 2).Value = 'US$2,581'
$ws.Cells.Item(14, 3).ClearContents()
$ws.Cells.Item(14, 4).ClearContents()
$ws.Cells.Item(14, 5).ClearContents()

$ws.Cells.Item(15, 1).Value = 'GuestReady - Modern Apt for 2 in Paris'
$ws.Cells.Item(15, 2).Value = 'US$3,102'
$ws.Cells.Item(15, 3).ClearContents()
$ws.Cells.Item(15, 4).ClearContents()
$ws.Cells.Item(15, 5).ClearContents()

$ws.Cells.Item(16, 1).Value = 'Hotel Regence Paris'
$ws.Cells.Item(16, 2).Value = 'US$4,451'
$ws.Cells.Item(16, 3).Value = '''6.7'
$ws.Cells.Item(16, 4).Value = 'Review score'
$ws.Cells.Item(16, 5).Value = '''370'

$ws.Cells.Item(17, 1).Value = 'Apartment hotel with view of Paris next to subway'
$ws.Cells.Item(17, 2).Value = 'US$2,622'
$ws.Cells.Item(17, 3).Value = '''7.6'
$ws.Cells.Item(17, 4).Value = 'Good'
$ws.Cells.Item(17, 5).Value = '''26'

$ws.Cells.Item(18, 1).Value = 'housewithgreenshutters'
$ws.Cells.Item(18, 2).Value = 'US$5,360'
$ws.Cells.Item(18, 3).ClearContents()
$ws.Cells.Item(18, 4).ClearContents()
$ws.Cells.Item(18, 5).ClearContents()

$ws.Cells.Item(19, 1).Value = 'Auteuil Cozy Studio'
$ws.Cells.Item(19, 2).Value = 'US$3,988'
$ws.Cells.Item(19, 3).Value = '''9.6'
$ws.Cells.Item(19, 4).Value = 'Exceptional'
$ws.Cells.Item(19, 5).Value = '''7'

$ws.Cells.Item(20, 1).Value = 'GuestReady - Artistic Gem in Central Paris'
$ws.Cells.Item(20, 2).Value = 'US$4,326'
$ws.Cells.Item(20, 3).ClearContents()
$ws.Cells.Item(20, 4).ClearContents()
$ws.Cells.Item(20, 5).ClearContents()

$ws.Cells.Item(21, 1).Value = 'GuestReady - White Bliss in the 20th'
$ws.Cells.Item(21, 2).Value = 'US$2,211'
$ws.Cells.Item(21, 3).ClearContents()
$ws.Cells.Item(21, 4).ClearContents()
$ws.Cells.Item(21, 5).ClearContents()

$ws.Cells.Item(22, 1).Value = 'GuestReady - My little comfort in the 10th'
$ws.Cells.Item(22, 2).Value = 'US$2,443'
$ws.Cells.Item(22, 3).Value = '''1.0'
$ws.Cells.Item(22, 4).Value = 'Review score'
$ws.Cells.Item(22, 5).Value = '''1'

$ws.Cells.Item(23, 1).Value = 'Studette cosy entièrement rénovée'
$ws.Cells.Item(23, 2).Value = 'US$2,621'
$ws.Cells.Item(23, 3).Value = '''6.6'
$ws.Cells.Item(23, 4).Value = 'Review score'
$ws.Cells.Item(23, 5).Value = '''11'

$ws.Cells.Item(24, 1).Value = 'PERE LACHAISE - MONTMARTRE Paris centre'
$ws.Cells.Item(24, 2).Value = 'US$3,458'
$ws.Cells.Item(24, 3).Value = '''6.6'
$ws.Cells.Item(24, 4).Value = 'Review score'
$ws.Cells.Item(24, 5).Value = '''50'

$ws.Cells.Item(25, 1).Value = 'Studio Montmartre!'
$ws.Cells.Item(25, 2).Value = 'US$1,638'
$ws.Cells.Item(25, 3).ClearContents()
$ws.Cells.Item(25, 4).ClearContents()
$ws.Cells.Item(25, 5).ClearContents()

$ws.Cells.Item(26, 1).Value = 'SAINT PAUL 22'
$ws.Cells.Item(26, 2).Value = 'US$3,331'
$ws.Cells.Item(26, 3).Value = '''5.9'
$ws.Cells.Item(26, 4).Value = 'Review score'
$ws.Cells.Item(26, 5).Value = '''37'

# Rows 27-28 no longer exist in the refreshed extract; clear them so the
# sheet's used range (and dimension) shrinks back down to A1:E26.
$ws.Range("A27:E28").ClearContents()

Write-Host "Edit complete"